# Update CDA Logical model for ST.r2b

$wb = $excel.ActiveWorkbook

# --- Rename the "Include from EntityNamePartQu" sheet to "Include #0" ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInc = $wb.Worksheets.Item("Include from EntityNamePartQu")
$wsInc.Name = "Include #0"

# --- Metadata sheet updates ---
# Bump version + date
$wsMeta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$wsMeta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row after "Contact" (row 10) for the new "Jurisdiction" property,
# pushing Description/Purpose/Copyright/Immutable down by one row.
$wsMeta.Rows.Item(11).Insert()

# Match the style of the surrounding table rows (Contact, row 10).
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""
